$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected for this product/market. It belongs
# chronologically right after the existing row 54, so insert a fresh row at
# position 55 (this pushes the former rows 55-61 down to 56-62, matching the
# data already present in those rows) and populate it with the new record.
$ws.Rows(55).Insert()

$ws.Cells.Item(55, 1).Value = 2
$ws.Cells.Item(55, 2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(55, 3).Value = 'Coquimbo'
$ws.Cells.Item(55, 4).Value = 44776
$ws.Cells.Item(55, 5).Value = 4
$ws.Cells.Item(55, 6).Value = 100112026
$ws.Cells.Item(55, 7).Value = 'Haba'
$ws.Cells.Item(55, 8).Value = 'Sin especificar'
$ws.Cells.Item(55, 9).Value = 'Primera'
$ws.Cells.Item(55, 10).Value = 1100
$ws.Cells.Item(55, 11).Value = 10000
$ws.Cells.Item(55, 12).Value = 11000
$ws.Cells.Item(55, 13).Value = 10500
$ws.Cells.Item(55, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(55, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(55, 16).Value = 420
$ws.Cells.Item(55, 17).Value = 25
$ws.Cells.Item(55, 18).Value = 'Hortaliza'
